$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the data rows (2-21) with the updated TPM values and the newly added
# "Inflammatory-Mac" cluster. Target-cluster order is ECs, FAPs, Inflammatory-Mac,
# MuSCs, Resolving-Mac for every sending cluster (ECs, FAPs, Inflammatory-Mac, MuSCs).

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna1"
$ws.Cells.Item(2, 3).Value = "Epha2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 14.02618833333333
$ws.Cells.Item(2, 8).Value = 42.078565
$ws.Cells.Item(2, 9).Value = 0.806325281849088
$ws.Cells.Item(2, 10).Value = 0.8172785134657441
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 12.997753
$ws.Cells.Item(2, 14).Value = 38.993259
$ws.Cells.Item(2, 15).Value = 0.4740421406233454
$ws.Cells.Item(2, 16).Value = 0.5546503645614554
$ws.Cells.Item(2, 17).Value = 182.3089314881483
$ws.Cells.Item(2, 18).Value = 1640.780383393335
$ws.Cells.Item(2, 19).Value = 0.382232162646464
$ws.Cells.Item(2, 20).Value = 0.4533038254420194

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna1"
$ws.Cells.Item(3, 3).Value = "Epha2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 14.02618833333333
$ws.Cells.Item(3, 8).Value = 42.078565
$ws.Cells.Item(3, 9).Value = 0.806325281849088
$ws.Cells.Item(3, 10).Value = 0.8172785134657441
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.7472513333333333
$ws.Cells.Item(3, 14).Value = 2.241754
$ws.Cells.Item(3, 15).Value = 0.02725306609819269
$ws.Cells.Item(3, 16).Value = 0.03188729809316786
$ws.Cells.Item(3, 17).Value = 10.48108793366778
$ws.Cells.Item(3, 18).Value = 94.32979140300999
$ws.Cells.Item(3, 19).Value = 0.02197483620287704
$ws.Cells.Item(3, 20).Value = 0.02606080358402329

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna1"
$ws.Cells.Item(4, 3).Value = "Epha2"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 14.02618833333333
$ws.Cells.Item(4, 8).Value = 42.078565
$ws.Cells.Item(4, 9).Value = 0.806325281849088
$ws.Cells.Item(4, 10).Value = 0.8172785134657441
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.182591666666666
$ws.Cells.Item(4, 14).Value = 3.547775
$ws.Cells.Item(4, 15).Value = 0.04313039993528083
$ws.Cells.Item(4, 16).Value = 0.05046448405689858
$ws.Cells.Item(4, 17).Value = 16.58725343809722
$ws.Cells.Item(4, 18).Value = 149.285280942875
$ws.Cells.Item(4, 19).Value = 0.0347771318840792
$ws.Cells.Item(4, 20).Value = 0.04124353851283782

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna1"
$ws.Cells.Item(5, 3).Value = "Epha2"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 14.02618833333333
$ws.Cells.Item(5, 8).Value = 42.078565
$ws.Cells.Item(5, 9).Value = 0.806325281849088
$ws.Cells.Item(5, 10).Value = 0.8172785134657441
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 11.9545335
$ws.Cells.Item(5, 14).Value = 23.909067
$ws.Cells.Item(5, 15).Value = 0.4359947946767024
$ws.Cells.Item(5, 16).Value = 0.3400888529957002
$ws.Cells.Item(5, 17).Value = 167.6765383081425
$ws.Cells.Item(5, 18).Value = 1006.059229848855
$ws.Cells.Item(5, 19).Value = 0.3515536257024274
$ws.Cells.Item(5, 20).Value = 0.2779473122225958

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Efna1"
$ws.Cells.Item(6, 3).Value = "Epha2"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 14.02618833333333
$ws.Cells.Item(6, 8).Value = 42.078565
$ws.Cells.Item(6, 9).Value = 0.806325281849088
$ws.Cells.Item(6, 10).Value = 0.8172785134657441
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.5368526666666666
$ws.Cells.Item(6, 14).Value = 1.610558
$ws.Cells.Item(6, 15).Value = 0.01957959866647858
$ws.Cells.Item(6, 16).Value = 0.022909000292778
$ws.Cells.Item(6, 17).Value = 7.529996609918888
$ws.Cells.Item(6, 18).Value = 67.76996948927
$ws.Cells.Item(6, 19).Value = 0.01578752541324037
$ws.Cells.Item(6, 20).Value = 0.0187230337042679

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna1"
$ws.Cells.Item(7, 3).Value = "Epha2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.483247333333333
$ws.Cells.Item(7, 8).Value = 7.449742
$ws.Cells.Item(7, 9).Value = 0.1427547569137158
$ws.Cells.Item(7, 10).Value = 0.1446939568272663
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 12.997753
$ws.Cells.Item(7, 14).Value = 38.993259
$ws.Cells.Item(7, 15).Value = 0.4740421406233454
$ws.Cells.Item(7, 16).Value = 0.5546503645614554
$ws.Cells.Item(7, 17).Value = 32.27663547657534
$ws.Cells.Item(7, 18).Value = 290.489719289178
$ws.Cells.Item(7, 19).Value = 0.06767177055154315
$ws.Cells.Item(7, 20).Value = 0.08025455590408276

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna1"
$ws.Cells.Item(8, 3).Value = "Epha2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.483247333333333
$ws.Cells.Item(8, 8).Value = 7.449742
$ws.Cells.Item(8, 9).Value = 0.1427547569137158
$ws.Cells.Item(8, 10).Value = 0.1446939568272663
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.7472513333333333
$ws.Cells.Item(8, 14).Value = 2.241754
$ws.Cells.Item(8, 15).Value = 0.02725306609819269
$ws.Cells.Item(8, 16).Value = 0.03188729809316786
$ws.Cells.Item(8, 17).Value = 1.855609880829778
$ws.Cells.Item(8, 18).Value = 16.700488927468
$ws.Cells.Item(8, 19).Value = 0.003890504826000926
$ws.Cells.Item(8, 20).Value = 0.004613899333631002

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna1"
$ws.Cells.Item(9, 3).Value = "Epha2"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.483247333333333
$ws.Cells.Item(9, 8).Value = 7.449742
$ws.Cells.Item(9, 9).Value = 0.1427547569137158
$ws.Cells.Item(9, 10).Value = 0.1446939568272663
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.182591666666666
$ws.Cells.Item(9, 14).Value = 3.547775
$ws.Cells.Item(9, 15).Value = 0.04313039993528083
$ws.Cells.Item(9, 16).Value = 0.05046448405689858
$ws.Cells.Item(9, 17).Value = 2.936667602672222
$ws.Cells.Item(9, 18).Value = 26.43000842405
$ws.Cells.Item(9, 19).Value = 0.006157069758352358
$ws.Cells.Item(9, 20).Value = 0.007301905877439153

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Efna1"
$ws.Cells.Item(10, 3).Value = "Epha2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.483247333333333
$ws.Cells.Item(10, 8).Value = 7.449742
$ws.Cells.Item(10, 9).Value = 0.1427547569137158
$ws.Cells.Item(10, 10).Value = 0.1446939568272663
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 11.9545335
$ws.Cells.Item(10, 14).Value = 23.909067
$ws.Cells.Item(10, 15).Value = 0.4359947946767024
$ws.Cells.Item(10, 16).Value = 0.3400888529957002
$ws.Cells.Item(10, 17).Value = 29.686063435119
$ws.Cells.Item(10, 18).Value = 178.116380610714
$ws.Cells.Item(10, 19).Value = 0.06224033092971808
$ws.Cells.Item(10, 20).Value = 0.04920880181279436

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Efna1"
$ws.Cells.Item(11, 3).Value = "Epha2"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.483247333333333
$ws.Cells.Item(11, 8).Value = 7.449742
$ws.Cells.Item(11, 9).Value = 0.1427547569137158
$ws.Cells.Item(11, 10).Value = 0.1446939568272663
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.5368526666666666
$ws.Cells.Item(11, 14).Value = 1.610558
$ws.Cells.Item(11, 15).Value = 0.01957959866647858
$ws.Cells.Item(11, 16).Value = 0.022909000292778
$ws.Cells.Item(11, 17).Value = 1.333137952892889
$ws.Cells.Item(11, 18).Value = 11.998241576036
$ws.Cells.Item(11, 19).Value = 0.002795080848101263
$ws.Cells.Item(11, 20).Value = 0.003314793899319052

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Efna1"
$ws.Cells.Item(12, 3).Value = "Epha2"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1863673333333333
$ws.Cells.Item(12, 8).Value = 0.559102
$ws.Cells.Item(12, 9).Value = 0.01071372271683668
$ws.Cells.Item(12, 10).Value = 0.01085925937435662
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 12.997753
$ws.Cells.Item(12, 14).Value = 38.993259
$ws.Cells.Item(12, 15).Value = 0.4740421406233454
$ws.Cells.Item(12, 16).Value = 0.5546503645614554
$ws.Cells.Item(12, 17).Value = 2.422356565935333
$ws.Cells.Item(12, 18).Value = 21.801209093418
$ws.Cells.Item(12, 19).Value = 0.005078756050734224
$ws.Cells.Item(12, 20).Value = 0.006023092170854304

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Efna1"
$ws.Cells.Item(13, 3).Value = "Epha2"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1863673333333333
$ws.Cells.Item(13, 8).Value = 0.559102
$ws.Cells.Item(13, 9).Value = 0.01071372271683668
$ws.Cells.Item(13, 10).Value = 0.01085925937435662
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.7472513333333333
$ws.Cells.Item(13, 14).Value = 2.241754
$ws.Cells.Item(13, 15).Value = 0.02725306609819269
$ws.Cells.Item(13, 16).Value = 0.03188729809316786
$ws.Cells.Item(13, 17).Value = 0.1392632383231111
$ws.Cells.Item(13, 18).Value = 1.253369144908
$ws.Cells.Item(13, 19).Value = 0.0002919817933596586
$ws.Cells.Item(13, 20).Value = 0.0003462724407411372

# Row 14
$ws.Cells.Item(14, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 2).Value = "Efna1"
$ws.Cells.Item(14, 3).Value = "Epha2"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.1863673333333333
$ws.Cells.Item(14, 8).Value = 0.559102
$ws.Cells.Item(14, 9).Value = 0.01071372271683668
$ws.Cells.Item(14, 10).Value = 0.01085925937435662
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.182591666666666
$ws.Cells.Item(14, 14).Value = 3.547775
$ws.Cells.Item(14, 15).Value = 0.04313039993528083
$ws.Cells.Item(14, 16).Value = 0.05046448405689858
$ws.Cells.Item(14, 17).Value = 0.2203964553388889
$ws.Cells.Item(14, 18).Value = 1.98356809805
$ws.Cells.Item(14, 19).Value = 0.0004620871455728695
$ws.Cells.Item(14, 20).Value = 0.0005480069215669463

# Row 15
$ws.Cells.Item(15, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 2).Value = "Efna1"
$ws.Cells.Item(15, 3).Value = "Epha2"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.1863673333333333
$ws.Cells.Item(15, 8).Value = 0.559102
$ws.Cells.Item(15, 9).Value = 0.01071372271683668
$ws.Cells.Item(15, 10).Value = 0.01085925937435662
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 11.9545335
$ws.Cells.Item(15, 14).Value = 23.909067
$ws.Cells.Item(15, 15).Value = 0.4359947946767024
$ws.Cells.Item(15, 16).Value = 0.3400888529957002
$ws.Cells.Item(15, 17).Value = 2.227934529639
$ws.Cells.Item(15, 18).Value = 13.367607177834
$ws.Cells.Item(15, 19).Value = 0.004671127336150331
$ws.Cells.Item(15, 20).Value = 0.003693113065007749

# Row 16
$ws.Cells.Item(16, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 2).Value = "Efna1"
$ws.Cells.Item(16, 3).Value = "Epha2"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.1863673333333333
$ws.Cells.Item(16, 8).Value = 0.559102
$ws.Cells.Item(16, 9).Value = 0.01071372271683668
$ws.Cells.Item(16, 10).Value = 0.01085925937435662
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.5368526666666666
$ws.Cells.Item(16, 14).Value = 1.610558
$ws.Cells.Item(16, 15).Value = 0.01957959866647858
$ws.Cells.Item(16, 16).Value = 0.022909000292778
$ws.Cells.Item(16, 17).Value = 0.1000517998795556
$ws.Cells.Item(16, 18).Value = 0.900466198916
$ws.Cells.Item(16, 19).Value = 0.0002097703910195967
$ws.Cells.Item(16, 20).Value = 0.0002487747761864881

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Efna1"
$ws.Cells.Item(17, 3).Value = "Epha2"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.6993955000000001
$ws.Cells.Item(17, 8).Value = 1.398791
$ws.Cells.Item(17, 9).Value = 0.04020623852035952
$ws.Cells.Item(17, 10).Value = 0.02716827033263282
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 12.997753
$ws.Cells.Item(17, 14).Value = 38.993259
$ws.Cells.Item(17, 15).Value = 0.4740421406233454
$ws.Cells.Item(17, 16).Value = 0.5546503645614554
$ws.Cells.Item(17, 17).Value = 9.090569958311502
$ws.Cells.Item(17, 18).Value = 54.54341974986901
$ws.Cells.Item(17, 19).Value = 0.01905945137460403
$ws.Cells.Item(17, 20).Value = 0.01506889104449897

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Efna1"
$ws.Cells.Item(18, 3).Value = "Epha2"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.6993955000000001
$ws.Cells.Item(18, 8).Value = 1.398791
$ws.Cells.Item(18, 9).Value = 0.04020623852035952
$ws.Cells.Item(18, 10).Value = 0.02716827033263282
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.7472513333333333
$ws.Cells.Item(18, 14).Value = 2.241754
$ws.Cells.Item(18, 15).Value = 0.02725306609819269
$ws.Cells.Item(18, 16).Value = 0.03188729809316786
$ws.Cells.Item(18, 17).Value = 0.5226242199023333
$ws.Cells.Item(18, 18).Value = 3.135745319414
$ws.Cells.Item(18, 19).Value = 0.001095743275955059
$ws.Cells.Item(18, 20).Value = 0.0008663227347724316

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Efna1"
$ws.Cells.Item(19, 3).Value = "Epha2"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.6993955000000001
$ws.Cells.Item(19, 8).Value = 1.398791
$ws.Cells.Item(19, 9).Value = 0.04020623852035952
$ws.Cells.Item(19, 10).Value = 0.02716827033263282
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.182591666666666
$ws.Cells.Item(19, 14).Value = 3.547775
$ws.Cells.Item(19, 15).Value = 0.04313039993528083
$ws.Cells.Item(19, 16).Value = 0.05046448405689858
$ws.Cells.Item(19, 17).Value = 0.8270992900041666
$ws.Cells.Item(19, 18).Value = 4.962595740025
$ws.Cells.Item(19, 19).Value = 0.0017341111472764
$ws.Cells.Item(19, 20).Value = 0.00137103274505466

# Row 20
$ws.Cells.Item(20, 1).Value = "MuSCs"
$ws.Cells.Item(20, 2).Value = "Efna1"
$ws.Cells.Item(20, 3).Value = "Epha2"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.6993955000000001
$ws.Cells.Item(20, 8).Value = 1.398791
$ws.Cells.Item(20, 9).Value = 0.04020623852035952
$ws.Cells.Item(20, 10).Value = 0.02716827033263282
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 11.9545335
$ws.Cells.Item(20, 14).Value = 23.909067
$ws.Cells.Item(20, 15).Value = 0.4359947946767024
$ws.Cells.Item(20, 16).Value = 0.3400888529957002
$ws.Cells.Item(20, 17).Value = 8.360946934499252
$ws.Cells.Item(20, 18).Value = 33.44378773799701
$ws.Cells.Item(20, 19).Value = 0.01752971070840667
$ws.Cells.Item(20, 20).Value = 0.009239625895302207

# Row 21
$ws.Cells.Item(21, 1).Value = "MuSCs"
$ws.Cells.Item(21, 2).Value = "Efna1"
$ws.Cells.Item(21, 3).Value = "Epha2"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.6993955000000001
$ws.Cells.Item(21, 8).Value = 1.398791
$ws.Cells.Item(21, 9).Value = 0.04020623852035952
$ws.Cells.Item(21, 10).Value = 0.02716827033263282
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.5368526666666666
$ws.Cells.Item(21, 14).Value = 1.610558
$ws.Cells.Item(21, 15).Value = 0.01957959866647858
$ws.Cells.Item(21, 16).Value = 0.022909000292778
$ws.Cells.Item(21, 17).Value = 0.3754723392296667
$ws.Cells.Item(21, 18).Value = 2.252834035378
$ws.Cells.Item(21, 19).Value = 0.0007872220141173508
$ws.Cells.Item(21, 20).Value = 0.0006223979130045573
